$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update recalculated values in columns D/E for existing rows
$ws.Range("D2").Value = 0.001016323201858174
$ws.Range("D3").Value = 0.004050965523330885
$ws.Range("D4").Value = 0.0149849239216109
$ws.Range("E4").Value = 0.8805219312989998
$ws.Range("D5").Value = 0.01180130576484925

# Add new row 6 for DWA, copying the formatting from row 5's label cell
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "DWA"

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
